# Adds two new experiment rows (60 and 61) to the "Experiments" log sheet:
#   Row 60 -> EXP58 (TPR-LSTM -> csv, run_id 48, from run 48)
#   Row 61 -> EXP59 (TPR concatenated with LSTM, run_id 52, 60,000 updates)
# Mirrors the existing rows (52 and 56-59) in layout/formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Write cell values in the exact order new unique strings are first
#    introduced, so the shared-string table is built up in the same order
#    as the target workbook (A60, E60, D60, B60, D61, A61, B61, ...).
# ---------------------------------------------------------------------------

$descA60 = 'aR & aF vector to excel for the whole test set, nRoles=20, TPR-LSTM, cR=0.00005, cF=0.00005. [running from "QA_TPR_for_Run" branch "master" ]'
$ws.Range("A60").Value = $descA60

$ws.Range("E60").Value = "from run_id=48"

$ws.Range("D60").Value = "EXP58.txt"

$cmdB60 = 'python -m basic.cli --mode test --LSTMandTPR True --TPRregularizer1 True --TPRvis True --write2csv True --which_tensors2vis "fw_u_aR,bw_u_aR,fw_u_aF,bw_u_aF" --nRoles 20 --vis True --batch_size 40 --load_path "/home/hpalangi/QA/TPR_Stuff/Codes/TPR_ver1.0/QA_TPR_for_Run/out/basic/48/save/basic-20000" --run_id 48 |& tee /home/hpalangi/QA/TPR_Stuff/Codes/TPR_ver1.0/Log_Files/EXP58.txt'
$ws.Range("B60").Value = $cmdB60

$ws.Range("C60").Value = "DLDGX / 1"
$ws.Range("F60").Value = 2

$ws.Range("D61").Value = "EXP59.txt"

$descA61 = 'TPR concatenated with LSTM in 
phrase embedding layer 
batchsize = 40. With visualizations. With regularization. Regularization weights=0.00005. num_steps=60,000,  nRoles=20, nSymbols=100 [running from "QA_TPR_for_Run" branch "master" ]. '
$ws.Range("A61").Value = $descA61

$cmdB61 = 'python -m basic.cli --mode train --noload --len_opt --cluster --LSTMandTPR True --TPRregularizer1 True --TPRvis True num_steps=60000 --cF 0.00005 --cR 0.00005 --nRoles 20 --nSymbols 100 --batch_size 40 --run_id 52 |& tee /home/hpalangi/QA/TPR_Stuff/Codes/TPR_ver1.0/Log_Files/EXP59.txt'
$ws.Range("B61").Value = $cmdB61

$ws.Range("C61").Value = "DLDGX / 3"
$ws.Range("E61").Value = 52
$ws.Range("F61").Value = 3

# ---------------------------------------------------------------------------
# 2. Copy cell-level formatting from the most similar existing rows so the
#    new rows pick up the same fill/border/number-format/alignment styles
#    already defined in the workbook (row 52 -> row 60, row 59 -> row 61).
# ---------------------------------------------------------------------------

foreach ($col in @("A","B","C","D","E","F","K","L")) {
    $src = $ws.Range($col + "52")
    $dst = $ws.Range($col + "60")
    $src.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats
}

foreach ($col in @("A","B","C","D","E","F","H","I","K","L")) {
    $src = $ws.Range($col + "59")
    $dst = $ws.Range($col + "61")
    $src.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Row heights matching the taller wrapped-text rows.
# ---------------------------------------------------------------------------

$ws.Rows.Item(60).RowHeight = 105
$ws.Rows.Item(61).RowHeight = 195

# ---------------------------------------------------------------------------
# 4. Rich-text (bold) runs inside the two description cells.
# ---------------------------------------------------------------------------

# A60: "aR & aF vector to excel for the whole test set, nRoles=20, " +
#      "TPR-LSTM"(bold) + ","(plain) + " cR=0.00005, cF=0.00005. "(bold) +
#      "[running from "QA_TPR_for_Run" branch "master" ]"(plain)
$ws.Range("A60").Characters(60, 8).Font.Bold = $true
$ws.Range("A60").Characters(68, 1).Font.Bold = $false
$ws.Range("A60").Characters(69, 25).Font.Bold = $true
$ws.Range("A60").Characters(94, 48).Font.Bold = $false

# A61: "TPR concatenated with LSTM"(bold) +
#      " in \nphrase embedding layer \nbatchsize = 40. With visualizations. With regularization. "(plain) +
#      "Regularization weights=0.00005. num_steps=60,000,  nRoles=20, nSymbols=100"(bold) +
#      " [running from "QA_TPR_for_Run" branch "master" ]. "(plain)
$ws.Range("A61").Characters(1, 26).Font.Bold = $true
$ws.Range("A61").Characters(27, 87).Font.Bold = $false
$ws.Range("A61").Characters(114, 74).Font.Bold = $true
$ws.Range("A61").Characters(188, 51).Font.Bold = $false

# ---------------------------------------------------------------------------
# 5. Move the active selection to B61 (matches the saved view of the file
#    after the edit).
# ---------------------------------------------------------------------------

$ws.Range("B61").Select()

Write-Output "done"
